# Clients.xlsx fix:
# The old sheet had an extra "CustomerID" column in column A, which pushed
# "ItemName" (and the customer's First/Last Name) one column to the right.
# The web-scraping tool expects ItemName (and eventually a Total) to be the
# very first column so it can overwrite the sheet without disturbing the
# other columns. Remove the CustomerID column entirely and shift the rest
# of the data left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A (CustomerID). This shifts ItemName -> A, First Name -> B,
# Last Name -> C, and automatically drops the now-unused "CustomerID" shared
# string.
$ws.Columns.Item(1).Delete()

# Reflect where the cursor ended up after the column removal (next column
# after the new data, i.e. E1 as recorded by Excel).
$ws.Range("E1").Select()
